# Generate Report for Handoff
# ----------------------------------------------------------------------------
# Two files were newly handed off for localization since the last report:
#   960e020f-5a52-4b29-9fa0-298f93c2548a.md  -> 3974610d-0856-4f76-94fa-99bda592d0c0.md
#   c78b7663-4a9f-4b67-948f-bde31b7bfa51.md  -> ffff3c37e763-1625-4cf4-862a-3e14f39d38ee.md
# Their status flips from "Handed back: in sync with en-US" to "Ready for
# handoff", the latest-handoff xliff/file metadata is refreshed, and the
# (now stale) handback file/date columns are cleared back to "no handback
# yet".
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldUuid1 = "960e020f-5a52-4b29-9fa0-298f93c2548a"
$newUuid1 = "3974610d-0856-4f76-94fa-99bda592d0c0"
$oldUuid2 = "c78b7663-4a9f-4b67-948f-bde31b7bfa51"
$newUuid2 = "ffff3c37e763-1625-4cf4-862a-3e14f39d38ee"

$newHash = "eb58e5290ab270ecee8ed2216865a42d49b4320d"

$newStatus = "Ready for handoff"
$overviewDate = "2016-08-16 00:57:32"
$zhHandoffDate = "2016-08-16 00:57:27"
$deHandoffDate = "2016-08-16 00:57:32"
$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newUuid1.md"
$ov.Range("A3").Value = "$newUuid2.md"

$ov.Range("B2").Value = "e2e\$newUuid1.md"
$ov.Range("B3").Value = "e2e\$newUuid2.md"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newUuid1.md"
    }
    elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newUuid2.md"
    }
}

$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("G2").Value = $overviewDate

$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Range("G3").Value = $overviewDate

$ov.Columns.Item(5).ColumnWidth = 16.25
$ov.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newUuid1.md"
$zh.Range("A3").Value = "$newUuid2.md"

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    }
    elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    }
    elseif ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $hl.Delete()
    }
}

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$zh.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"

$zh.Range("H2").Value = $zhHandoffDate
$zh.Range("H3").Value = $zhHandoffDate

$zh.Range("F3").Value = "True"

$zh.Range("I2").Value = ""
$zh.Range("I3").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("I3").Style = "Normal"

$zh.Range("J2").Value = ""
$zh.Range("J3").Value = ""

$zh.Range("K2").Value = $noHandback
$zh.Range("K3").Value = $noHandback

$zh.Columns.Item(3).ColumnWidth = 16.25
$zh.Columns.Item(9).ColumnWidth = 17.75
$zh.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newUuid1.md"
$de.Range("A3").Value = "$newUuid2.md"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    }
    elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    }
    elseif ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $hl.Delete()
    }
}

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$de.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"

$de.Range("H2").Value = $deHandoffDate
$de.Range("H3").Value = $deHandoffDate

$de.Range("F3").Value = "True"

$de.Range("I2").Value = ""
$de.Range("I3").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("I3").Style = "Normal"

$de.Range("J2").Value = ""
$de.Range("J3").Value = ""

$de.Range("K2").Value = $noHandback
$de.Range("K3").Value = $noHandback

$de.Columns.Item(3).ColumnWidth = 16.25
$de.Columns.Item(9).ColumnWidth = 17.75
$de.Columns.Item(10).ColumnWidth = 20.75
